# "fix(gui) step 1 and 2"
# Bumps the quote date by one day and updates the two support-pair prices
# on the "SOP. CAÑOS CURVOS" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1: quote date, 45308 (2024-01-17) -> 45309 (2024-01-18)
$ws.Range("A1").Value = 45309

# D30: price for "Par Soporte p/ CAÑO CURVO NATURAL" (SOCC-010)
$ws.Range("D30").Value = 1475

# D31: price for "Par Soporte p/ CAÑO CURVO BLANCO" (SOCC-011)
$ws.Range("D31").Value = 1680
